$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 31950
$ws.Range("J3").Value = 31950
$ws.Range("L3").Value = 31950
$ws.Range("N3").Value = -32178
$ws.Range("H17").Value = 1115.7354
$ws.Range("I17").Value = 750
$ws.Range("J17").Value = 1151.129
$ws.Range("K17").Value = 2250
$ws.Range("L17").Value = 3453.387
$ws.Range("M17").Value = -2082
$ws.Range("N17").Value = -3789.387
$ws.Range("H102").Value = 31950
$ws.Range("J102").Value = 31950
$ws.Range("L102").Value = 31950
$ws.Range("N102").Value = -38440
$ws.Range("H113").Value = 3151.7097
$ws.Range("I113").Value = 2255.4
$ws.Range("J113").Value = 3578.524
$ws.Range("K113").Value = 2255.4
$ws.Range("L113").Value = 3578.524
$ws.Range("M113").Value = 998.5999999999999
$ws.Range("N113").Value = -10086.524
$ws.Range("H116").Value = 39695.43
$ws.Range("I116").Value = 54333.3
$ws.Range("K116").Value = 54333.3
$ws.Range("M116").Value = -50891.3
$ws.Range("H132").Value = 4916.775
$ws.Range("I132").Value = 2468.2693
$ws.Range("J132").Value = 9464
$ws.Range("K132").Value = 7404.8079
$ws.Range("L132").Value = 28392
$ws.Range("M132").Value = -4874.8079
$ws.Range("N132").Value = -33452

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 738.86957
$ws.Range("I2").Value = 699.7
$ws.Range("K2").Value = 699.7
$ws.Range("M2").Value = -586.7
$ws.Range("H74").Value = 2099.6135
$ws.Range("I74").Value = 1199.4242
$ws.Range("J74").Value = 4800.1816
$ws.Range("K74").Value = 1199.4242
$ws.Range("L74").Value = 4800.1816
$ws.Range("M74").Value = -325.4241999999999
$ws.Range("N74").Value = -6548.1816
$ws.Range("H77").Value = 2099.6135
$ws.Range("I77").Value = 1199.4242
$ws.Range("J77").Value = 4800.1816
$ws.Range("K77").Value = 5997.120999999999
$ws.Range("L77").Value = 24000.908
$ws.Range("M77").Value = -1629.120999999999
$ws.Range("N77").Value = -32736.908
$ws.Range("H92").Value = 28300
$ws.Range("J92").Value = 28300
$ws.Range("L92").Value = 28300
$ws.Range("N92").Value = -33292
$ws.Range("H116").Value = 738.86957
$ws.Range("I116").Value = 699.7
$ws.Range("K116").Value = 699.7
$ws.Range("M116").Value = 1594.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 738.86957
$ws.Range("I3").Value = 699.7
$ws.Range("K3").Value = 699.7
$ws.Range("M3").Value = -585.7
$ws.Range("H75").Value = 2957
$ws.Range("I75").Value = 2957
$ws.Range("K75").Value = 2957
$ws.Range("M75").Value = -2021
$ws.Range("H78").Value = 2957
$ws.Range("I78").Value = 2957
$ws.Range("K78").Value = 8871
$ws.Range("M78").Value = -4191
$ws.Range("H86").Value = 5318.5356
$ws.Range("I86").Value = 4454.722
$ws.Range("J86").Value = 6873.4
$ws.Range("K86").Value = 4454.722
$ws.Range("L86").Value = 6873.4
$ws.Range("M86").Value = -3331.722
$ws.Range("N86").Value = -9119.4
$ws.Range("H89").Value = 5318.5356
$ws.Range("I89").Value = 4454.722
$ws.Range("J89").Value = 6873.4
$ws.Range("K89").Value = 22273.61
$ws.Range("L89").Value = 34367
$ws.Range("M89").Value = -16657.61
$ws.Range("N89").Value = -45599

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 10029.214
$ws.Range("I69").Value = 6540.9
$ws.Range("J69").Value = 18750
$ws.Range("K69").Value = 6540.9
$ws.Range("L69").Value = 18750
$ws.Range("M69").Value = -5791.9
$ws.Range("N69").Value = -20248
$ws.Range("H72").Value = 10029.214
$ws.Range("I72").Value = 6540.9
$ws.Range("J72").Value = 18750
$ws.Range("K72").Value = 19622.7
$ws.Range("L72").Value = 56250
$ws.Range("M72").Value = -15878.7
$ws.Range("N72").Value = -63738
$ws.Range("H81").Value = 20000
$ws.Range("J81").Value = 20000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -21996
$ws.Range("H82").Value = 28933.334
$ws.Range("I82").Value = 9000
$ws.Range("J82").Value = 38900
$ws.Range("K82").Value = 9000
$ws.Range("L82").Value = 38900
$ws.Range("M82").Value = -8639
$ws.Range("N82").Value = -39622
$ws.Range("H84").Value = 20000
$ws.Range("J84").Value = 20000
$ws.Range("L84").Value = 60000
$ws.Range("N84").Value = -69984
$ws.Range("H85").Value = 28933.334
$ws.Range("I85").Value = 9000
$ws.Range("J85").Value = 38900
$ws.Range("K85").Value = 9000
$ws.Range("L85").Value = 38900
$ws.Range("M85").Value = -7752
$ws.Range("N85").Value = -41396
$ws.Range("H132").Value = 1801.283
$ws.Range("I132").Value = 1015.68
$ws.Range("J132").Value = 2502.7144
$ws.Range("K132").Value = 3047.04
$ws.Range("L132").Value = 7508.1432
$ws.Range("M132").Value = -517.04
$ws.Range("N132").Value = -12568.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 212.5
$ws.Range("I14").Value = 212.5
$ws.Range("K14").Value = 637.5
$ws.Range("M14").Value = -464.5
$ws.Range("H70").Value = 3964.6
$ws.Range("I70").Value = 911.5
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 2734.5
$ws.Range("L70").Value = 18000
$ws.Range("M70").Value = -2419.5
$ws.Range("N70").Value = -18630
$ws.Range("H73").Value = 3964.6
$ws.Range("I73").Value = 911.5
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 2734.5
$ws.Range("L73").Value = 18000
$ws.Range("M73").Value = -1642.5
$ws.Range("N73").Value = -20184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3641.75
$ws.Range("I80").Value = 3772.9412
$ws.Range("K80").Value = 3772.9412
$ws.Range("M80").Value = -2774.9412
$ws.Range("H83").Value = 3641.75
$ws.Range("I83").Value = 3772.9412
$ws.Range("K83").Value = 18864.706
$ws.Range("M83").Value = -13872.706
$ws.Range("H87").Value = 35500
$ws.Range("J87").Value = 35500
$ws.Range("L87").Value = 35500
$ws.Range("N87").Value = -37996
$ws.Range("H90").Value = 35500
$ws.Range("J90").Value = 35500
$ws.Range("L90").Value = 106500
$ws.Range("N90").Value = -118980
$ws.Range("H107").Value = 5237.3335
$ws.Range("I107").Value = 10258
$ws.Range("J107").Value = 673.0909
$ws.Range("K107").Value = 10258
$ws.Range("L107").Value = 673.0909
$ws.Range("M107").Value = -8338
$ws.Range("N107").Value = -4513.0909

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 42581.1
$ws.Range("J46").Value = 42581.1
$ws.Range("L46").Value = 42581.1
$ws.Range("N46").Value = -43043.1
$ws.Range("H132").Value = 19712.232
$ws.Range("I132").Value = 29712.4
$ws.Range("J132").Value = 3045.2856
$ws.Range("K132").Value = 89137.20000000001
$ws.Range("L132").Value = 9135.856800000001
$ws.Range("M132").Value = -86607.20000000001
$ws.Range("N132").Value = -14195.8568
$ws.Range("H134").Value = 42581.1
$ws.Range("J134").Value = 42581.1
$ws.Range("L134").Value = 127743.3
$ws.Range("N134").Value = -132813.3
$ws.Range("H136").Value = 45979576
$ws.Range("I136").Value = 90911910
$ws.Range("J136").Value = 18520924
$ws.Range("K136").Value = 272735730
$ws.Range("L136").Value = 55562772
$ws.Range("M136").Value = -272733180
$ws.Range("N136").Value = -55567872
